$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(401, 9, 48, 67, 75, 45)
    3 = @(1202, 2, 10, 10, 10, 10)
    4 = @(1001, 18, 30, 75, 60, 72)
    5 = @(601, 9, 60, 67, 60, 42)
    6 = @(1203, 3, 15, 15, 15, 15)
    7 = @(301, 6, 45, 30, 60, 45)
    8 = @(501, 9, 52, 30, 75, 45)
    9 = @(201, 9, 30, 15, 45, 30)
    10 = @(701, 3, 90, 45, 97, 15)
    11 = @(801, 3, 67, 65, 52, 45)
    12 = @(1201, 2, 10, 10, 10, 10)
    13 = @(101, 9, 30, 15, 60, 15)
    14 = @(901, 16, 15, 45, 60, 60)
    15 = @(902, 1, 0, 0, 0, 0)
    16 = @(1, 0, 2, 2, 2, 2)
    17 = @(2, 0, 2, 2, 2, 2)
    18 = @(802, 0, 4, 5, 4, 0)
    19 = @(1101, 0, 15, 30, 30, 0)
    20 = @(3, 0, 3, 3, 3, 3)
    21 = @(502, 0, 4, 0, 0, 0)
    22 = @(402, 0, 0, 4, 0, 0)
    23 = @(602, 0, 0, 4, 0, 9)
    24 = @(702, 0, 0, 0, 4, 0)
    25 = @(1002, 0, 0, 0, 0, 9)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item([int]$r, $i + 1).Value = $vals[$i]
    }
}
